{"js": "// Replace the date and each two-digit multiplication problem's text with\n// its updated value, as described by the diff. Every old value is unique\n// in the document, so a simple exact search-and-replace per pair is safe.\nconst replacements = [\n  [\"2024-03-14 Thursday\", \"2024-03-15 Friday\"],\n  [\"89\u00d780=\", \"19\u00d752=\"],\n  [\"78\u00d752=\", \"30\u00d794=\"],\n  [\"21\u00d753=\", \"53\u00d750=\"],\n  [\"50\u00d753=\", \"85\u00d796=\"],\n  [\"85\u00d765=\", \"18\u00d771=\"],\n  [\"48\u00d777=\", \"31\u00d748=\"],\n  [\"81\u00d719=\", \"30\u00d740=\"],\n  [\"69\u00d722=\", \"21\u00d797=\"],\n  [\"15\u00d757=\", \"47\u00d788=\"],\n  [\"49\u00d730=\", \"49\u00d781=\"],\n  [\"76\u00d722=\", \"51\u00d739=\"],\n  [\"47\u00d791=\", \"99\u00d727=\"],\n  [\"25\u00d766=\", \"89\u00d759=\"],\n  [\"59\u00d750=\", \"13\u00d780=\"],\n  [\"71\u00d786=\", \"25\u00d738=\"],\n  [\"52\u00d746=\", \"95\u00d781=\"],\n  [\"64\u00d754=\", \"39\u00d773=\"],\n  [\"37\u00d726=\", \"63\u00d775=\"],\n  [\"30\u00d746=\", \"29\u00d730=\"],\n  [\"11\u00d722=\", \"79\u00d750=\"],\n  [\"82\u00d751=\", \"40\u00d735=\"],\n  [\"49\u00d768=\", \"47\u00d764=\"],\n  [\"61\u00d760=\", \"82\u00d799=\"],\n  [\"65\u00d790=\", \"52\u00d759=\"],\n  [\"91\u00d759=\", \"24\u00d726=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each two-digit multiplication problem's text with\n# its updated value, as described by the diff. Every old value is unique\n# in the document, so a plain Find/Replace (MatchWildcards off) per pair\n# is safe and will touch only the one intended run.\n$d = $word.ActiveDocument\n\n$olds = @(\n  \"2024-03-14 Thursday\",\n  \"89\u00d780=\",\n  \"78\u00d752=\",\n  \"21\u00d753=\",\n  \"50\u00d753=\",\n  \"85\u00d765=\",\n  \"48\u00d777=\",\n  \"81\u00d719=\",\n  \"69\u00d722=\",\n  \"15\u00d757=\",\n  \"49\u00d730=\",\n  \"76\u00d722=\",\n  \"47\u00d791=\",\n  \"25\u00d766=\",\n  \"59\u00d750=\",\n  \"71\u00d786=\",\n  \"52\u00d746=\",\n  \"64\u00d754=\",\n  \"37\u00d726=\",\n  \"30\u00d746=\",\n  \"11\u00d722=\",\n  \"82\u00d751=\",\n  \"49\u00d768=\",\n  \"61\u00d760=\",\n  \"65\u00d790=\",\n  \"91\u00d759=\"\n)\n\n$news = @(\n  \"2024-03-15 Friday\",\n  \"19\u00d752=\",\n  \"30\u00d794=\",\n  \"53\u00d750=\",\n  \"85\u00d796=\",\n  \"18\u00d771=\",\n  \"31\u00d748=\",\n  \"30\u00d740=\",\n  \"21\u00d797=\",\n  \"47\u00d788=\",\n  \"49\u00d781=\",\n  \"51\u00d739=\",\n  \"99\u00d727=\",\n  \"89\u00d759=\",\n  \"13\u00d780=\",\n  \"25\u00d738=\",\n  \"95\u00d781=\",\n  \"39\u00d773=\",\n  \"63\u00d775=\",\n  \"29\u00d730=\",\n  \"79\u00d750=\",\n  \"40\u00d735=\",\n  \"47\u00d764=\",\n  \"82\u00d799=\",\n  \"52\u00d759=\",\n  \"24\u00d726=\"\n)\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n  $find = $d.Content.Find\n  $find.Text = $olds[$i]\n  $find.Replacement.Text = $news[$i]\n  $find.Execute($olds[$i], $false, $false, $false, $false, $false, $true, 1, $false, $news[$i], 2)\n}\n"}
